$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (column F) values for matching events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 467
$wsExhibit.Range("F3").Value = 5587
$wsExhibit.Range("F5").Value = 67
$wsExhibit.Range("F8").Value = 52
$wsExhibit.Range("F9").Value = 529
$wsExhibit.Range("F10").Value = 22

# Sheet "全部类型" - same events appear again (different row numbers) - update F column
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 467
$wsAll.Range("F3").Value = 5587
$wsAll.Range("F6").Value = 67
$wsAll.Range("F10").Value = 52
$wsAll.Range("F11").Value = 529
$wsAll.Range("F12").Value = 22
